$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("Elapsed") values are stored as native numbers (t="n").
$numericUpdates = @{
    "B2" = 12
    "B3" = 1
    "B4" = 7
    "B5" = 8
    "B6" = 486
    "B7" = 8
    "B8" = 26
    "B9" = 2
    "B10" = 7
    "B11" = 581
    "B12" = 11
    "B14" = 62
    "B15" = 29
    "B17" = 17027
    "B20" = 62
    "B21" = 2
    "B22" = 38
    "B23" = 8
    "B24" = 46
    "B25" = 7
    "B26" = 5
    "B27" = 11
    "B28" = 10380
    "B29" = 419
    "B30" = 35
    "B31" = 6
    "B32" = 6
    "B33" = 444
    "B34" = 10
    "B35" = 37
    "B36" = 134
    "B37" = 5
    "B38" = 0
    "B39" = 49
    "B40" = 47
    "B41" = 49
    "B42" = 4
    "B43" = 2
    "B44" = 48
    "B45" = 45
    "B46" = 5
    "B47" = 61
    "B48" = 7
    "B49" = 11
    "B50" = 9
    "B51" = 6
    "B52" = 6
    "B53" = 2
    "B54" = 50
    "B55" = 505
    "B56" = 48
    "B57" = 6
    "B58" = 494
    "B59" = 1
    "B60" = 5
    "B61" = 44
}

foreach ($ref in $numericUpdates.Keys) {
    $ws.Range($ref).Value = $numericUpdates[$ref]
}

# Columns C-L are stored as text (inline strings) even though their
# contents look numeric, so each new value must be written as a string.
# Forcing the Text number format before the write keeps Excel from
# re-interpreting the string as a number, and resetting the style back
# to "Normal" afterwards avoids leaving a stray number-format override.
$textUpdates = @{
    "C6" = "11198"
    "D6" = "31"
    "E6" = "3497"
    "F6" = "68"
    "G6" = "3429"
    "H6" = "3497"
    "I6" = "341384"
    "J6" = "7433"
    "K6" = "3765"
    "L6" = "7440"
    "C8" = "247"
    "D8" = "47"
    "E8" = "269"
    "F8" = "23"
    "G8" = "246"
    "H8" = "269"
    "I8" = "15629"
    "J8" = "212"
    "L8" = "213"
    "C11" = "5296"
    "D11" = "57"
    "E11" = "5192"
    "F11" = "515"
    "G11" = "4677"
    "H11" = "5192"
    "I11" = "337799"
    "J11" = "4193"
    "K11" = "1103"
    "L11" = "4194"
    "E14" = "457"
    "G14" = "452"
    "H14" = "457"
    "C15" = "268"
    "D15" = "46"
    "E15" = "325"
    "F15" = "29"
    "G15" = "296"
    "H15" = "325"
    "I15" = "16936"
    "J15" = "233"
    "K15" = "35"
    "L15" = "234"
    "E17" = "180807"
    "G17" = "180750"
    "H17" = "180807"
    "E20" = "457"
    "G20" = "452"
    "H20" = "457"
    "D23" = "71"
    "E23" = "133"
    "F23" = "32"
    "G23" = "101"
    "H23" = "133"
    "C28" = "158970"
    "D28" = "196"
    "E28" = "29757"
    "G28" = "28621"
    "H28" = "29757"
    "I28" = "9275719"
    "K28" = "158933"
    "E30" = "302"
    "G30" = "299"
    "H30" = "302"
    "C35" = "209"
    "D35" = "160"
    "E35" = "411"
    "F35" = "13"
    "G35" = "398"
    "H35" = "411"
    "I35" = "11836"
    "K35" = "193"
    "L35" = "18"
    "C36" = "978"
    "D36" = "45"
    "E36" = "867"
    "F36" = "77"
    "G36" = "790"
    "H36" = "867"
    "I36" = "70553"
    "J36" = "848"
    "K36" = "130"
    "L36" = "849"
    "E42" = "42"
    "G42" = "40"
    "H42" = "42"
    "C45" = "424"
    "D45" = "60"
    "E45" = "342"
    "F45" = "60"
    "G45" = "282"
    "H45" = "342"
    "I45" = "27092"
    "J45" = "418"
    "K45" = "6"
    "L45" = "370"
    "E46" = "69"
    "G46" = "67"
    "H46" = "69"
    "E47" = "281"
    "G47" = "279"
    "H47" = "281"
    "C50" = "56"
    "D50" = "40"
    "E50" = "42"
    "G50" = "38"
    "H50" = "42"
    "I50" = "4116"
    "J50" = "44"
    "L50" = "45"
    "D54" = "29"
    "E54" = "410"
    "G54" = "404"
    "H54" = "410"
    "C55" = "11312"
    "E55" = "3625"
    "F55" = "78"
    "G55" = "3547"
    "H55" = "3625"
    "I55" = "345826"
    "J55" = "7505"
    "K55" = "3807"
    "L55" = "7514"
    "C58" = "5379"
    "E58" = "4096"
    "F58" = "33"
    "G58" = "4063"
    "H58" = "4096"
    "I58" = "294670"
    "J58" = "5027"
    "K58" = "352"
    "L58" = "5032"
    "D61" = "90"
    "E61" = "128"
    "F61" = "52"
    "G61" = "76"
    "H61" = "128"
}

foreach ($ref in $textUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$ref]
    $cell.Style = "Normal"
}
